# Nhập excel vào kho đến
# The template's "Xe nâng" column (column D) is removed entirely; every
# column to its right shifts one place to the left (E->D, F->E, ... Z->Y).
# This also relocates the per-row comments that were anchored on column T
# (they now live on column S) and drops the now-unused "Xe nâng" shared
# string plus the trailing column Z formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Comments don't automatically re-anchor when a whole column is deleted via
# COM, so capture their text up front, remove them, delete the column, and
# then re-create them one cell to the left (T -> S) at the same rows.
$commentRows = 4, 5, 6, 7, 8, 9
$commentText = @{}
foreach ($r in $commentRows) {
    $cell = $ws.Range("T$r")
    if ($cell.Comment -ne $null) {
        $commentText[$r] = $cell.Comment.Text()
        $cell.Comment.Delete()
    }
}

# Delete column D ("Xe nâng") -- shifts everything from E onward left by one.
$ws.Columns.Item(4).Delete()

foreach ($r in $commentRows) {
    if ($commentText.ContainsKey($r)) {
        $ws.Range("S$r").AddComment($commentText[$r]) | Out-Null
    }
}

# Restore the active selection to match the post-edit workbook state.
$ws.Range("G20").Select()
